$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Career Projection sheet: bump the projected years forward by one, and
# become the active tab/selection.
# ---------------------------------------------------------------------------
$wsProj = $wb.Worksheets.Item("Career Projection")
$wsProj.Range("A2").Value = 2022
$wsProj.Range("A3").Value = 2023
$wsProj.Range("A4").Value = 2024
$wsProj.Range("A5").Value = 2025

# ---------------------------------------------------------------------------
# Assets sheet: restructure from (Account Name/Type, Account Balance,
# Earning estimate) into (Priority, Type, Balance, Growth Percent), add a
# Priority column, rename the USAA account to Checking, reorder rows, and
# apply currency formatting to the Balance column.
# ---------------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Range("A1:C5").ClearContents()

# Introduce brand-new shared strings in this exact order so the saved
# sharedStrings.xml table matches the canonical append order.
$wsAssets.Range("A1").Value = "Priority"
$wsAssets.Range("D1").Value = "Growth Percent"
$wsAssets.Range("C1").Value = "Balance"
$wsAssets.Range("B1").Value = "Type"
$wsAssets.Range("A1:D1").Font.Bold = $true

$wsAssets.Range("A2").Value = 1
$wsAssets.Range("B2").Value = "TSP"
$wsAssets.Range("C2").Value = 10000
$wsAssets.Range("D2").Value = 6

$wsAssets.Range("A3").Value = 2
$wsAssets.Range("B3").Value = "Roth IRA"
$wsAssets.Range("C3").Value = 15000
$wsAssets.Range("D3").Value = 6

$wsAssets.Range("A4").Value = 3
$wsAssets.Range("B4").Value = "Savings"
$wsAssets.Range("C4").Value = 7500
$wsAssets.Range("D4").Value = 1.45

$wsAssets.Range("A5").Value = 4
$wsAssets.Range("B5").Value = "Checking"
$wsAssets.Range("C5").Value = 5000
$wsAssets.Range("D5").Value = 0

$wsAssets.Range("C2:C5").NumberFormat = '"$"#,##0_);[Red]("$"#,##0)'
$wsAssets.PageSetup.Orientation = 1

# Narrower columns to fit the new Priority/Type/Balance layout (closest the
# engine's pixel-grid column-width rounding can reach to the authored sizes).
$wsAssets.Columns.Item(1).ColumnWidth = 5.75
$wsAssets.Columns.Item(2).ColumnWidth = 8.6
$wsAssets.Columns.Item(3).ColumnWidth = 9.25

# ---------------------------------------------------------------------------
# Window / tab-activation state: Career Projection becomes the active sheet
# (previously it was Assets).
# ---------------------------------------------------------------------------
$wsAssets.Activate()
$wsAssets.Range("L10").Select()
$wsProj.Activate()
$wsProj.Range("A6").Select()
